$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows (358-366), continuing the existing series through 1/09/2021
$data = @(
    @(358, 44432, 1, 6, 68.99724011039559),
    @(359, 44433, 2, 8, 91.99632014719411),
    @(360, 44434, 0, 7, 80.49678012879485),
    @(361, 44435, 1, 8, 91.99632014719411),
    @(362, 44436, 5, 11, 126.4949402023919),
    @(363, 44437, 0, 9, 103.4958601655934),
    @(364, 44438, 3, 12, 137.9944802207912),
    @(365, 44439, 0, 11, 126.4949402023919),
    @(366, 44440, 0, 9, 103.4958601655934)
)

# Column A uses the same date/time style (border + bold + centered + custom
# date number format) as the rows immediately above; replicate it onto the
# new rows by copying the formatting from the last existing row.
$ws.Range("A357").Copy()
$ws.Range("A358:A366").PasteSpecial(-4122)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
